$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")
$ws.Activate()

# CUDA (column I) timings added for SUM ARRAY (row 6), MIN ARRAY (row 9) and MAT X VEC (row 10).
# D6/D9/D10 were placeholders (text) that now hold the real measured millisecond values.
$ws.Range("D6").Value = 215.29519999999999
$ws.Range("I6").Value = 73.250900000000001

$ws.Range("D9").Value = 0.81159999999999999
$ws.Range("I9").Value = 6.6848999999999998

$ws.Range("D10").Value = 29615.030599999998
$ws.Range("I10").Value = 10041.270399999999

# Speedup column (I) for CUDA, mirroring the existing H-column formulas (C/H).
# I20/I27 are brand-new cells with the sheet's default (unstyled) format; I23/I24
# already existed as empty, styled placeholder cells and must keep their style.
$ws.Range("I20").Formula = "=D6/I6"
$ws.Range("I20").Style = "Normal"

$ws.Range("I23").Formula = "=D9/I9"
$ws.Range("I24").Formula = "=D10/I10"

$ws.Range("I27").Formula = "=D13/I13"
$ws.Range("I27").Style = "Normal"

$excel.Calculate()

# Restore the view/selection recorded at save time.
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
